$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.584.25"
$ws.Range("E2").Value = "  -7.41%  "
$ws.Range("D3").Value = "1.684.87"
$ws.Range("E3").Value = "  -6.57%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.75"
$ws.Range("E5").Value = "  -6.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4958"
$ws.Range("E7").Value = "  -16.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2596"
$ws.Range("E8").Value = "  -6.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.59"
$ws.Range("E9").Value = "  -7.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06077"
$ws.Range("E10").Value = "  -11.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07268"
$ws.Range("E11").Value = "  -3.75%  "
$ws.Range("D12").Value = "1.665.42"
$ws.Range("E12").Value = "  -7.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.423"
$ws.Range("E13").Value = "  -6.48%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5709"
$ws.Range("E14").Value = "  -9.11%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "1.914.54"
$ws.Range("E15").Value = "  -6.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008199"
$ws.Range("E16").Value = "  -12.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.41"
$ws.Range("E17").Value = "  -14.68%  "
$ws.Range("D18").Value = "26.604.47"
$ws.Range("E18").Value = "  -7.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.986"
$ws.Range("E19").Value = "  -9.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.68"
$ws.Range("E21").Value = "  -6.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "182.03"
$ws.Range("E22").Value = "  -14.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.172"
$ws.Range("E23").Value = "  -10.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("E25").Value = "  -6.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.536"
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1128"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.19"
$ws.Range("E28").Value = "  -7.56%  "
$ws.Range("E29").Value = "  -9.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05581"
$ws.Range("E30").Value = "  -10.73%  "
$ws.Range("E31").Value = "  -7.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.465"
$ws.Range("E32").Value = "  -8.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.447"
$ws.Range("E33").Value = "  -8.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.636"
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("E35").Value = "  -5.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.394"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5834"
$ws.Range("E37").Value = "  -9.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.606"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01569"
$ws.Range("E39").Value = "  -8.32%  "
$ws.Range("D40").Value = "1.069.69"
$ws.Range("E40").Value = "  -6.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.863"
$ws.Range("E41").Value = "  -8.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8464"
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.11"
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("D45").Value = "1.841.09"
$ws.Range("E45").Value = "  -6.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.14"
$ws.Range("E46").Value = "  -7.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000104"
$ws.Range("E47").Value = "  -6.88%  "
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.068"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4334"
$ws.Range("E50").Value = "  -3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05176"
$ws.Range("E51").Value = "  -5.26%  "
